$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Page break paragraph, right after the "About meetings..." para.
# ------------------------------------------------------------------
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$pageBreakRange = $d.Paragraphs.Last.Range
$pageBreakXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:r><w:br w:type="page"/></w:r></w:p>'
$pageBreakRange.InsertXML($pageBreakXml)

# ------------------------------------------------------------------
# 2) "Notes for next meeting" heading (bold, 14pt / sz 28) split
#    across three runs, first run carries lastRenderedPageBreak.
# ------------------------------------------------------------------
$afterBreak = $d.Paragraphs.Last.Range
$afterBreak.Collapse(0)
$afterBreak.InsertParagraphAfter()
$headingRange = $d.Paragraphs.Last.Range
$headingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:lastRenderedPageBreak/><w:t>Note</w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>s</w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> for next meeting</w:t></w:r>' +
  '</w:p>'
$headingRange.InsertXML($headingXml)

# ------------------------------------------------------------------
# 3) New bulleted (List Paragraph / numId 1) question paragraph.
# ------------------------------------------------------------------
$afterHeading = $d.Paragraphs.Last.Range
$afterHeading.Collapse(0)
$afterHeading.InsertParagraphAfter()
$bulletRange = $d.Paragraphs.Last.Range
$bulletXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t>How should we filter the navigation status in the pre-processing? Which ones can be considered as relevant to detect anomalies and which ones are just noise?</w:t></w:r>' +
  '</w:p>'
$bulletRange.InsertXML($bulletXml)

# ------------------------------------------------------------------
# 4) Normal style now carries an explicit en-GB language.
# ------------------------------------------------------------------
$normalStyle = $d.Styles("Normal")
$normalStyle.LanguageID = "en-GB"

Write-Output "done"
